# Make the header row's column names dynamic "keys" style (uppercase),
# matching the keys now used from the data object: CRNO, NAME, AGE, SEX,
# TOPICS, and a new WORK column replacing the previously-blank F1 header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CRNO"
$ws.Range("B1").Value = "NAME"
$ws.Range("C1").Value = "AGE"
$ws.Range("D1").Value = "SEX"
$ws.Range("E1").Value = "TOPICS"
$ws.Range("F1").Value = "WORK"
